{"js": "// M2Doc #253: refresh the recorded Java stack-trace text (POI 3.16 -> 3.17)\n// shown in the \"conditionRuntimeException\" sample so line numbers / frames\n// match the new dependency versions.\n\nconst startAnchor = \"divOp(java.lang.Integer\";\nconst endAnchor = \"RemoteTestRunner.main(RemoteTestRunner.java:192)\";\n\nconst startResults = context.document.body.search(startAnchor, { matchCase: true });\nconst endResults = context.document.body.search(endAnchor, { matchCase: true });\nawait context.sync();\n\nif (startResults.items.length !== 1 || endResults.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly one match for each anchor, got \" +\n    startResults.items.length + \" / \" + endResults.items.length\n  );\n}\n\n// The whole stack trace lives in a single bold/red run; replace it in one\n// shot (much cheaper & safer than chasing individual line diffs) with the\n// updated trace text, keeping the trailing \"\\n\" (the run already ends with\n// one, right before the manual line break) untouched.\nconst newStackTraceLines = [\n  \"divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:\",\n  \"\\t/ by zero\",\n  \"java.lang.ArithmeticException: / by zero\",\n  \"\\tat org.eclipse.acceleo.query.services.NumberServices.divOp(NumberServices.java:99)\",\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\",\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\",\n  \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\",\n  \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\",\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:163)\",\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:136)\",\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)\",\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:168)\",\n  \"\\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)\",\n  \"\\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\",\n  \"\\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:183)\",\n  \"\\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\",\n  \"\\tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)\",\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseConditional(M2DocEvaluator.java:1267)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseConditional(M2DocEvaluator.java:1)\",\n  \"\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:134)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)\",\n  \"\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)\",\n  \"\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:264)\",\n  \"\\tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:712)\",\n  \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)\",\n  \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)\",\n  \"\\tat sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)\",\n  \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\",\n  \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\",\n  \"\\tat org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:50)\",\n  \"\\tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)\",\n  \"\\tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)\",\n  \"\\tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)\",\n  \"\\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\",\n  \"\\tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)\",\n  \"\\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)\",\n  \"\\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)\",\n  \"\\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\",\n  \"\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\",\n  \"\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\",\n  \"\\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\",\n  \"\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\",\n  \"\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\",\n  \"\\tat org.junit.runners.Suite.runChild(Suite.java:128)\",\n  \"\\tat org.junit.runners.Suite.runChild(Suite.java:27)\",\n  \"\\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\",\n  \"\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\",\n  \"\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\",\n  \"\\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\",\n  \"\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\",\n  \"\\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\",\n  \"\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\",\n  \"\\tat org.junit.runners.Suite.runChild(Suite.java:128)\",\n  \"\\tat org.junit.runners.Suite.runChild(Suite.java:27)\",\n  \"\\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\",\n  \"\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\",\n  \"\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\",\n  \"\\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\",\n  \"\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\",\n  \"\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\",\n  \"\\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\",\n  \"\\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\",\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:539)\",\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:761)\",\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:461)\",\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:207)\"\n];\n\nconst fullRange = startResults.items[0].expandTo(endResults.items[0]);\nfullRange.insertText(newStackTraceLines.join(\"\\n\"), Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# M2Doc #253: refresh the recorded Java stack-trace text (POI 3.16 -> 3.17)\n# shown in the \"conditionRuntimeException\" sample so line numbers / frames\n# match the new dependency versions.\n\n$d = $word.ActiveDocument\n\n# Locate the single bold/red run holding the stack trace via two unique\n# text anchors (its first and last line) and build a Range spanning the\n# whole block, rather than chasing each individual line diff.\n$startRange = $d.Content\n$startRange.Find.ClearFormatting()\n$startRange.Find.Text = \"divOp(java.lang.Integer\"\n$startRange.Find.Forward = $true\n$startRange.Find.Wrap = 0\n$startRange.Find.MatchCase = $true\n$startRange.Find.MatchWildcards = $false\n$found1 = $startRange.Find.Execute()\n\n$endRange = $d.Content\n$endRange.Find.ClearFormatting()\n$endRange.Find.Text = \"RemoteTestRunner.main(RemoteTestRunner.java:192)\"\n$endRange.Find.Forward = $true\n$endRange.Find.Wrap = 0\n$endRange.Find.MatchCase = $true\n$endRange.Find.MatchWildcards = $false\n$found2 = $endRange.Find.Execute()\n\nif (-not $found1 -or -not $found2) {\n    throw \"Could not locate the stack trace text to replace (found1=$found1 found2=$found2)\"\n}\n\n$target = $d.Range($startRange.Start, $endRange.End)\n\n# Updated stack trace text. The run's final newline character (right before\n# the manual line break) is left untouched, so this here-string intentionally\n# carries no trailing newline of its own.\n$newText = @'\ndivOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:\n\t/ by zero\njava.lang.ArithmeticException: / by zero\n\tat org.eclipse.acceleo.query.services.NumberServices.divOp(NumberServices.java:99)\n\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\n\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\n\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\n\tat java.lang.reflect.Method.invoke(Method.java:498)\n\tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:163)\n\tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:136)\n\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)\n\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:168)\n\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)\n\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:183)\n\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)\n\tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseConditional(M2DocEvaluator.java:1267)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseConditional(M2DocEvaluator.java:1)\n\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:134)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)\n\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)\n\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:264)\n\tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:712)\n\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)\n\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)\n\tat sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)\n\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\n\tat java.lang.reflect.Method.invoke(Method.java:498)\n\tat org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:50)\n\tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)\n\tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)\n\tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)\n\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\n\tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)\n\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)\n\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)\n\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\n\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\n\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\n\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\n\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\n\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\n\tat org.junit.runners.Suite.runChild(Suite.java:128)\n\tat org.junit.runners.Suite.runChild(Suite.java:27)\n\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\n\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\n\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\n\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\n\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\n\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\n\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\n\tat org.junit.runners.Suite.runChild(Suite.java:128)\n\tat org.junit.runners.Suite.runChild(Suite.java:27)\n\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\n\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\n\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\n\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\n\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\n\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\n\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\n\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:539)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:761)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:461)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:207)\n'@\n\n$target.Text = $newText\n"}
